$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Marno"
$ws.Range("A4").Value = "Phillip"
$ws.Range("A5").Value = "Abdel"

$ws.Range("B9").Select() | Out-Null
